# Implements "clear dtc info": adds ClearDtcInfo request and
# ClearDtcInfoResp response rows to the ReqResp table on the "ReqResp" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ReqResp")

# New rows appended after the existing data (which ends at row 77).
$newRows = @(
    @("ClearDtcInfo",     "sid",             "Req",  1, 1),
    @("ClearDtcInfo",     "dtcHighByte",     "Req",  2, 1),
    @("ClearDtcInfo",     "dtcMiddleByte",   "Req",  3, 1),
    @("ClearDtcInfo",     "dtcLowByte",      "Req",  4, 1),
    @("ClearDtcInfo",     "memorySelection", "Req",  5, 1),
    @("ClearDtcInfoResp", "sid",             "Resp", 1, 1)
)

$startRow = 78
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Resize the structured table ("Table8") to cover the new rows.
$table = $ws.ListObjects.Item("Table8")
$table.Resize($ws.Range("A1:F83"))

# Update the view to match where the author ended up after adding the rows
# (zoom level and final selection; the frozen header row stays as-is).
$ws.Activate()
$window = $excel.ActiveWindow
$window.Zoom = 90
$window.ScrollRow = 62
$window.ScrollColumn = 1
$ws.Range("B85").Select()
